$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E1").Value = "Insurance"
$ws.Range("E2").Value = 30
$ws.Range("E3").Value = -20
$ws.Range("E4").Value = 4.12
$ws.Range("E5").Value = -8
$ws.Range("E6").Value = -12

$ws.Range("E1:E6").Select()
